$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the three test-breakdown cells on the previous last data
# row (109) that no longer carry values in the updated export ---
$ws.Range("X109").ClearContents()
$ws.Range("Z109").ClearContents()
$ws.Range("AB109").ClearContents()

# --- Insert a new data row for 2020-06-16 (serial 43998) above the
# "TOTAL" footer row, which shifts that footer from row 110 to row 111
# and copies row 109's number formats/styles onto the new row ---
$ws.Rows.Item(110).Insert()

$ws.Range("A110").Value = 43998
$ws.Range("B110").Value = 2923
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 0.09
$ws.Range("E110").Value = 0.08
$ws.Range("F110").Value = 761
$ws.Range("G110").Value = 9
$ws.Range("H110").Value = 452
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = 40
$ws.Range("K110").Value = 6
$ws.Range("L110").Value = 4
$ws.Range("M110").Value = 0
$ws.Range("N110").Value = 0
$ws.Range("O110").Value = 0
$ws.Range("P110").Value = 2
$ws.Range("Q110").Value = 0
$ws.Range("R110").Value = 0
$ws.Range("S110").Value = 0
$ws.Range("T110").Value = 209
$ws.Range("U110").Value = 2668
$ws.Range("V110").Formula = "=B110-SUM(T110:U110)"
$ws.Range("W110").Value = 77584
$ws.Range("X110").Value = 13347
$ws.Range("Y110").Value = 53320
$ws.Range("Z110").Value = 9173
$ws.Range("AA110").Value = 24264
$ws.Range("AB110").Value = 4174
$ws.Range("AC110").Value = 2

# --- Keep the view anchored on the new last row, matching how Excel
# re-centers the selection after adding a row at the bottom of the
# table ---
$ws.Range("A111").Select() | Out-Null
